# Update countries & provincias Spain
# Refresh the COVID-19 country table on sheet "Pais": a handful of
# countries received new totals which shuffled their rank among the
# lowest-case countries (so several rows now show a different country
# name with updated figures), several existing countries had their
# daily figures revised, and the "last updated" timestamp advanced
# from 18:22 to 18:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 8 de Abril de 2020 a las 18:52"

# Row 4
$ws.Cells.Item(4, 2).Value = 406697
$ws.Cells.Item(4, 3).Value = 6362
$ws.Cells.Item(4, 4).Value = 22033
$ws.Cells.Item(4, 5).Value = 370796
$ws.Cells.Item(4, 7).Value = 1027
$ws.Cells.Item(4, 8).Value = 13868

# Row 17
$ws.Cells.Item(17, 2).Value = 14275
$ws.Cells.Item(17, 3).Value = 241
$ws.Cells.Item(17, 5).Value = 13442
$ws.Cells.Item(17, 7).Value = 20
$ws.Cells.Item(17, 8).Value = 706

# Row 26
$ws.Cells.Item(26, 4).Value = 506
$ws.Cells.Item(26, 5).Value = 5065
$ws.Cells.Item(26, 7).Value = 18
$ws.Cells.Item(26, 8).Value = 178

# Row 29
$ws.Cells.Item(29, 2).Value = 5402
$ws.Cells.Item(29, 3).Value = 331
$ws.Cells.Item(29, 5).Value = 3563

# Row 31
$ws.Cells.Item(31, 2).Value = 5205
$ws.Cells.Item(31, 3).Value = 357
$ws.Cells.Item(31, 5).Value = 4824
$ws.Cells.Item(31, 7).Value = 30
$ws.Cells.Item(31, 8).Value = 159

# Row 53
$ws.Cells.Item(53, 5).Value = 1314
$ws.Cells.Item(53, 7).Value = 3
$ws.Cells.Item(53, 8).Value = 63

# Row 112
$ws.Cells.Item(112, 4).Value = 50
$ws.Cells.Item(112, 5).Value = 155

# Row 133
$ws.Cells.Item(133, 2).Value = 81
$ws.Cells.Item(133, 3).Value = 2
$ws.Cells.Item(133, 5).Value = 76

# Row 140
$ws.Cells.Item(140, 2).Value = 59
$ws.Cells.Item(140, 3).Value = 3
$ws.Cells.Item(140, 4).Value = 16
$ws.Cells.Item(140, 5).Value = 36
$ws.Cells.Item(140, 7).Value = 2
$ws.Cells.Item(140, 8).Value = 7

# Row 158
$ws.Cells.Item(158, 1).Value = "Haiti"
$ws.Cells.Item(158, 2).Value = 27
$ws.Cells.Item(158, 3).Value = 2
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 5).Value = 26

# Row 159
$ws.Cells.Item(159, 1).Value = "Benin"
$ws.Cells.Item(159, 2).Value = 26
$ws.Cells.Item(159, 4).Value = 5
$ws.Cells.Item(159, 5).Value = 20

# Row 181
$ws.Cells.Item(181, 1).Value = "Somalia"
$ws.Cells.Item(181, 2).Value = 12
$ws.Cells.Item(181, 3).Value = 4
$ws.Cells.Item(181, 4).Value = 1
$ws.Cells.Item(181, 5).Value = 10
$ws.Cells.Item(181, 7).Value = 1
$ws.Cells.Item(181, 8).Value = 1

# Row 182
$ws.Cells.Item(182, 1).Value = "Suazilandia"
$ws.Cells.Item(182, 2).Value = 12
$ws.Cells.Item(182, 3).Value = 2
$ws.Cells.Item(182, 4).Value = 4
$ws.Cells.Item(182, 5).Value = 8

# Row 183
$ws.Cells.Item(183, 1).Value = "Seychelles"
$ws.Cells.Item(183, 5).Value = 11
$ws.Cells.Item(183, 8).Value = 0

# Row 184
$ws.Cells.Item(184, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 11

# Row 185
$ws.Cells.Item(185, 1).Value = "Zimbabue"
$ws.Cells.Item(185, 2).Value = 11
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 9
$ws.Cells.Item(185, 8).Value = 2

# Row 186
$ws.Cells.Item(186, 1).Value = "Groenlandia"
$ws.Cells.Item(186, 2).Value = 11
$ws.Cells.Item(186, 4).Value = 10
$ws.Cells.Item(186, 5).Value = 1

# Row 187
$ws.Cells.Item(187, 1).Value = "Republica del Chad"
$ws.Cells.Item(187, 4).Value = 2
$ws.Cells.Item(187, 5).Value = 8
$ws.Cells.Item(187, 8).Value = 0

# Row 188
$ws.Cells.Item(188, 1).Value = "Surinam"
$ws.Cells.Item(188, 2).Value = 10
$ws.Cells.Item(188, 4).Value = 3
$ws.Cells.Item(188, 5).Value = 6
$ws.Cells.Item(188, 8).Value = 1

# Row 189
$ws.Cells.Item(189, 1).Value = "Nepal"
$ws.Cells.Item(189, 4).Value = 1
$ws.Cells.Item(189, 5).Value = 8
$ws.Cells.Item(189, 8).Value = 0

# Row 190
$ws.Cells.Item(190, 1).Value = "Montserrat"
$ws.Cells.Item(190, 2).Value = 9
$ws.Cells.Item(190, 5).Value = 7
$ws.Cells.Item(190, 8).Value = 2

# Row 191
$ws.Cells.Item(191, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 8

# Row 192
$ws.Cells.Item(192, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(192, 4).Value = 1
$ws.Cells.Item(192, 8).Value = 0

# Row 195
$ws.Cells.Item(195, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 5).Value = 7
$ws.Cells.Item(195, 7).Value = 0

# Row 207
$ws.Cells.Item(207, 1).Value = "Burundi"

# Row 208
$ws.Cells.Item(208, 1).Value = "Anguila"

# Row 209
$ws.Cells.Item(209, 1).Value = "Islas Virgenes Britanicas"

# Row 210
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"

# Row 211
$ws.Cells.Item(211, 1).Value = "Papua Nueva Guinea"
